$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# J2 must stay a text string "001" (not numeric 1), so force text format first
$ws.Range("J2").NumberFormat = "@"
$ws.Range("J2").Value = "001"

$ws.Range("N2").Value = "2018-12-31 00:00:00"

$ws.Range("O2").Value = 365659823.42
$ws.Range("P2").Value = 7852453.7
$ws.Range("Q2").Value = 118238393
$ws.Range("R2").Value = 41.0652993612
$ws.Range("S2").Value = 151688292.44
$ws.Range("T2").Value = 41.647920473
$ws.Range("U2").Value = 2323863.35
$ws.Range("V2").Value = -0.4618059942
$ws.Range("W2").Value = 174759485.82
$ws.Range("X2").Value = 46490259.4
$ws.Range("Y2").Value = 57.4163296134
$ws.Range("Z2").Value = 6525680.49
$ws.Range("AA2").Value = 447.5476626407
$ws.Range("AB2").Value = 190900337.6
$ws.Range("AC2").Value = 46.4021093763
$ws.Range("AD2").Value = 29.4168838144
$ws.Range("AE2").Value = 14.8602879931
$ws.Range("AF2").Value = 168.8469555963
$ws.Range("AG2").Value = 47.7929142408
